$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Cfh", "Sell", "ECs", 3, 1, 0.7080493333333333, 2.124148, 0.008172603973273591, 0.008303215621510202, 2, 0.6666666666666666, 0.14839, 0.44517, 0.002798961202302754, 0.002799915128928837, 0.1050674405733333, 0.94560696516, 0.00002287480144297811, 0.00002324829903742467),
  @("ECs", "Cfh", "Sell", "Inflammatory-Mac", 3, 1, 0.7080493333333333, 2.124148, 0.008172603973273591, 0.008303215621510202, 3, 1, 50.610032, 151.830096, 0.9546163219576846, 0.9549416690637682, 35.83439941757867, 322.509594758208, 0.007801701145783195, 0.007929086584201306),
  @("ECs", "Cfh", "Sell", "MuSCs", 3, 1, 0.7080493333333333, 2.124148, 0.008172603973273591, 0.008303215621510202, 1, 0.5, 0.0541875, 0.108375, 0.001022095223059374, 0.0006816290453032835, 0.03836742325, 0.2302045395, 0.000008353179481038995, 0.000005659712937037308),
  @("ECs", "Cfh", "Sell", "Resolving-Mac", 3, 1, 0.7080493333333333, 2.124148, 0.008172603973273591, 0.008303215621510202, 2, 0.6666666666666666, 2.203488, 6.610464, 0.04156262161695323, 0.04157678676199977, 1.560178209408, 14.041603884672, 0.0003396748465663788, 0.000345221025334435),
  @("FAPs", "Cfh", "Sell", "ECs", 3, 1, 43.16235166666667, 129.487055, 0.4981980634967507, 0.5061600876489541, 2, 0.6666666666666666, 0.14839, 0.44517, 0.002798961202302754, 0.002799915128928837, 6.404861363816666, 57.64375227435, 0.001394437050789769, 0.001417205287068253),
  @("FAPs", "Cfh", "Sell", "Inflammatory-Mac", 3, 1, 43.16235166666667, 129.487055, 0.4981980634967507, 0.5061600876489541, 3, 1, 50.610032, 151.830096, 0.9546163219576846, 0.9549416690637682, 2184.447999045253, 19660.03199140728, 0.4755880029817092, 0.4833533589129555),
  @("FAPs", "Cfh", "Sell", "MuSCs", 3, 1, 43.16235166666667, 129.487055, 0.4981980634967507, 0.5061600876489541, 1, 0.5, 0.0541875, 0.108375, 0.001022095223059374, 0.0006816290453032835, 2.3388599309375, 14.033159585625, 0.0005092058608374595, 0.0003450134173147829),
  @("FAPs", "Cfh", "Sell", "Resolving-Mac", 3, 1, 43.16235166666667, 129.487055, 0.4981980634967507, 0.5061600876489541, 2, 0.6666666666666666, 2.203488, 6.610464, 0.04156262161695323, 0.04157678676199977, 95.10772394928, 855.9695155435201, 0.02070641760341429, 0.02104451003161568),
  @("Inflammatory-Mac", "Cfh", "Sell", "ECs", 3, 1, 20.47870633333333, 61.436119, 0.2363738638936219, 0.2401515068656985, 2, 0.6666666666666666, 0.14839, 0.44517, 0.002798961202302754, 0.002799915128928837, 3.038835232803333, 27.34951709523, 0.0006616012742766393, 0.0006724038373083265),
  @("Inflammatory-Mac", "Cfh", "Sell", "Inflammatory-Mac", 3, 1, 20.47870633333333, 61.436119, 0.2363738638936219, 0.2401515068656985, 3, 1, 50.610032, 151.830096, 0.9546163219576846, 0.9549416690637682, 1036.427982848603, 9327.851845637424, 0.2256463485570557, 0.2293306807945091),
  @("Inflammatory-Mac", "Cfh", "Sell", "MuSCs", 3, 1, 20.47870633333333, 61.436119, 0.2363738638936219, 0.2401515068656985, 1, 0.5, 0.0541875, 0.108375, 0.001022095223059374, 0.0006816290453032835, 1.1096898994375, 6.658139396625, 0.0002415965971417575, 0.000163694242353011),
  @("Inflammatory-Mac", "Cfh", "Sell", "Resolving-Mac", 3, 1, 20.47870633333333, 61.436119, 0.2363738638936219, 0.2401515068656985, 2, 0.6666666666666666, 2.203488, 6.610464, 0.04156262161695323, 0.04157678676199977, 45.124583661024, 406.121252949216, 0.00982431746514781, 0.009984727991528069),
  @("MuSCs", "Cfh", "Sell", "ECs", 2, 1, 4.0884615, 8.176923, 0.04719074664214938, 0.0319632882405021, 2, 0.6666666666666666, 0.14839, 0.44517, 0.002798961202302754, 0.002799915128928837, 0.606686801985, 3.64012081191, 0.0001320850689590751, 0.00008949449431489503),
  @("MuSCs", "Cfh", "Sell", "Inflammatory-Mac", 2, 1, 4.0884615, 8.176923, 0.04719074664214938, 0.0319632882405021, 3, 1, 50.610032, 151.830096, 0.9546163219576846, 0.9549416690637682, 206.917167345768, 1241.503004074608, 0.0450490569899656, 0.03052307582115139),
  @("MuSCs", "Cfh", "Sell", "MuSCs", 2, 1, 4.0884615, 8.176923, 0.04719074664214938, 0.0319632882405021, 1, 0.5, 0.0541875, 0.108375, 0.001022095223059374, 0.0006816290453032835, 0.22154350753125, 0.8861740301250001, 0.00004823343671554606, 0.00002178710564812712),
  @("MuSCs", "Cfh", "Sell", "Resolving-Mac", 2, 1, 4.0884615, 8.176923, 0.04719074664214938, 0.0319632882405021, 2, 0.6666666666666666, 2.203488, 6.610464, 0.04156262161695323, 0.04157678676199977, 9.008875853712, 54.05325512227201, 0.001961371146509161, 0.001328930819387691),
  @("Resolving-Mac", "Cfh", "Sell", "ECs", 3, 1, 18.199363, 54.598089, 0.2100647219942044, 0.2134219016233352, 2, 0.6666666666666666, 0.14839, 0.44517, 0.002798961202302754, 0.002799915128928837, 2.70060347557, 24.30543128013, 0.000587963006834292, 0.0005975632111999382),
  @("Resolving-Mac", "Cfh", "Sell", "Inflammatory-Mac", 3, 1, 18.199363, 54.598089, 0.2100647219942044, 0.2134219016233352, 3, 1, 50.610032, 151.830096, 0.9546163219576846, 0.9549416690637682, 921.070343809616, 8289.633094286544, 0.2005312122831709, 0.2038054669509511),
  @("Resolving-Mac", "Cfh", "Sell", "MuSCs", 3, 1, 18.199363, 54.598089, 0.2100647219942044, 0.2134219016233352, 1, 0.5, 0.0541875, 0.108375, 0.001022095223059374, 0.0006816290453032835, 0.9861779825625001, 5.917067895375, 0.0002147061488835717, 0.0001454745670503253),
  @("Resolving-Mac", "Cfh", "Sell", "Resolving-Mac", 3, 1, 18.199363, 54.598089, 0.2100647219942044, 0.2134219016233352, 2, 0.6666666666666666, 2.203488, 6.610464, 0.04156262161695323, 0.04157678676199977, 40.102077978144, 360.918701803296, 0.008730840555315589, 0.008873396894133901)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowData[$j]
    }
}
